$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dummy")

# Update the username/password test values in row 2, keep Results as "Pass"
$ws.Range("B2").Value = "fhsfh"
$ws.Range("C2").Value = "asfahs"
$ws.Range("D2").Value = "Pass"

# Move the active selection to H4 (matches the saved cursor position)
$ws.Activate()
$ws.Range("H4").Select()
